$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C
$ws.Columns.Item(3).ColumnWidth = 24.139196

# New rows of data
$ws.Range("A11").Value = "pear"
$ws.Range("B11").Value = "30"
$ws.Range("C11").Value = "btnOrchidItemTheme"
$ws.Range("D11").Value = "DarkGreen"

$ws.Range("A21").Value = "wet"
$ws.Range("B21").Value = "11"
$ws.Range("C21").Value = "btnDefaultItemTheme"
$ws.Range("D21").Value = "Green"

$ws.Range("A34").Value = "test"
$ws.Range("B34").Value = "32"
$ws.Range("C34").Value = "btnLimeGreenItemTheme"
$ws.Range("D34").Value = "Red"
